# EPI-551: Validate survey question calculation and writeToField.fieldType
#
# Add sample `calculation` values to the CalculatedQuestion / Result rows in
# both the "Question Validation Succeed" and "Question Validation Fail"
# sheets, and fix the bogus `fieldType` value in the PatientData "full
# config" sample JSON on the "Succeed" sheet so it uses a real field type.

$wb = $excel.ActiveWorkbook

$succeed = $wb.Worksheets.Item("Question Validation Succeed")
$fail = $wb.Worksheets.Item("Question Validation Fail")

# -- "Question Validation Fail" sheet ---------------------------------
# CalculatedQuestion rows (6-9) and Result rows (10-13): give each pair of
# validationCriteria/config rows a sample calculation.
$fail.Range("O6").Value = "1+1"
$fail.Range("O7").Value = "2+2"
$fail.Range("O8").Value = "1+1"
$fail.Range("O9").Value = "2+2"
$fail.Range("O10").Value = "1+1"
$fail.Range("O11").Value = "2+2"
$fail.Range("O12").Value = "1+1"
$fail.Range("O13").Value = "2+2"

# -- "Question Validation Succeed" sheet -------------------------------
# Succeed-CalculatedQuestion / Succeed-Result rows get a sample calculation.
$succeed.Range("O3").Value = "1+1"
$succeed.Range("O4").Value = "2+2"

# Succeed-PatientData "full config" sample used a made-up "def" fieldType;
# fix it to a real writeToPatient fieldType value.
$succeed.Range("P5").Value = '{ "column": "xyz", "writeToPatient": { "fieldName": "abc", "isAdditionalData": false, "fieldType": "FreeText" } }'

# -- Selection / active sheet bookkeeping ------------------------------
# Update the remembered selection on the Fail sheet first (leaves it as a
# background sheet), then on the Succeed sheet last so it stays the active
# tab, matching the workbook's activeTab.
$fail.Range("F14").Select()
$succeed.Range("P6").Select()
